$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1570.1428
$ws.Range("I9").Value = 1158.2
$ws.Range("K9").Value = 1158.2
$ws.Range("M9").Value = -989.2
$ws.Range("H17").Value = 2712.7144
$ws.Range("J17").Value = 2712.7144
$ws.Range("L17").Value = 8138.1432
$ws.Range("N17").Value = -8474.143199999999
$ws.Range("H29").Value = 5749.3335
$ws.Range("J29").Value = 5749.3335
$ws.Range("L29").Value = 17248.0005
$ws.Range("N29").Value = -17810.0005
$ws.Range("H38").Value = 5834.5
$ws.Range("I38").Value = 1252
$ws.Range("J38").Value = 14999.5
$ws.Range("K38").Value = 3756
$ws.Range("L38").Value = 44998.5
$ws.Range("M38").Value = -3384
$ws.Range("N38").Value = -45742.5
$ws.Range("H58").Value = 35
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H64").Value = 4150
$ws.Range("J64").Value = 4150
$ws.Range("L64").Value = 4150
$ws.Range("N64").Value = -4646
$ws.Range("H67").Value = 4150
$ws.Range("J67").Value = 4150
$ws.Range("L67").Value = 4150
$ws.Range("N67").Value = -5866
$ws.Range("H70").Value = 3837.3333
$ws.Range("I70").Value = 1283
$ws.Range("J70").Value = 4475.9165
$ws.Range("K70").Value = 3849
$ws.Range("L70").Value = 13427.7495
$ws.Range("M70").Value = -3579
$ws.Range("N70").Value = -13967.7495
$ws.Range("H73").Value = 3837.3333
$ws.Range("I73").Value = 1283
$ws.Range("J73").Value = 4475.9165
$ws.Range("K73").Value = 3849
$ws.Range("L73").Value = 13427.7495
$ws.Range("M73").Value = -2913
$ws.Range("N73").Value = -15299.7495
$ws.Range("H80").Value = 13845
$ws.Range("J80").Value = 30099.2
$ws.Range("L80").Value = 90297.60000000001
$ws.Range("N80").Value = -92293.60000000001
$ws.Range("H83").Value = 13845
$ws.Range("J83").Value = 30099.2
$ws.Range("L83").Value = 270892.8
$ws.Range("N83").Value = -280876.8
$ws.Range("H100").Value = 700.2727
$ws.Range("I100").Value = 570.3
$ws.Range("K100").Value = 570.3
$ws.Range("M100").Value = -29.29999999999995
$ws.Range("H111").Value = 1494.3
$ws.Range("I111").Value = 1201.375
$ws.Range("K111").Value = 3604.125
$ws.Range("M111").Value = -537.125
$ws.Range("H137").Value = 3123.2812
$ws.Range("I137").Value = 1291.7368
$ws.Range("K137").Value = 3875.2104
$ws.Range("M137").Value = -1325.2104
$ws.Range("H138").Value = 20645.75
$ws.Range("I138").Value = 25062
$ws.Range("K138").Value = 75186
$ws.Range("M138").Value = -70046

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 8416.666999999999
$ws.Range("I63").Value = 7500
$ws.Range("J63").Value = 9333.333000000001
$ws.Range("K63").Value = 7500
$ws.Range("L63").Value = 9333.333000000001
$ws.Range("M63").Value = -6814
$ws.Range("N63").Value = -10705.333
$ws.Range("H66").Value = 8416.666999999999
$ws.Range("I66").Value = 7500
$ws.Range("J66").Value = 9333.333000000001
$ws.Range("K66").Value = 37500
$ws.Range("L66").Value = 46666.665
$ws.Range("M66").Value = -34068
$ws.Range("N66").Value = -53530.665
$ws.Range("H132").Value = 8251.5
$ws.Range("I132").Value = 7377.75
$ws.Range("J132").Value = 9999
$ws.Range("K132").Value = 22133.25
$ws.Range("L132").Value = 29997
$ws.Range("M132").Value = -19603.25
$ws.Range("N132").Value = -35057

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2461.1667
$ws.Range("J86").Value = 3056
$ws.Range("L86").Value = 3056
$ws.Range("N86").Value = -5302
$ws.Range("H89").Value = 2461.1667
$ws.Range("J89").Value = 3056
$ws.Range("L89").Value = 15280
$ws.Range("N89").Value = -26512
$ws.Range("H105").Value = 4921.727
$ws.Range("I105").Value = 3870.5625
$ws.Range("J105").Value = 7724.8335
$ws.Range("K105").Value = 3870.5625
$ws.Range("L105").Value = 7724.8335
$ws.Range("M105").Value = -2123.5625
$ws.Range("N105").Value = -11218.8335
$ws.Range("H107").Value = 2828.6667
$ws.Range("I107").Value = 2203.1
$ws.Range("J107").Value = 4079.8
$ws.Range("K107").Value = 2203.1
$ws.Range("L107").Value = 4079.8
$ws.Range("M107").Value = -283.0999999999999
$ws.Range("N107").Value = -7919.8
$ws.Range("H134").Value = 5923.25
$ws.Range("I134").Value = 5850
$ws.Range("K134").Value = 17550
$ws.Range("M134").Value = -15015

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 549
$ws.Range("J22").Value = 313
$ws.Range("L22").Value = 313
$ws.Range("N22").Value = -1013
$ws.Range("H99").Value = 9687.111000000001
$ws.Range("I99").Value = 6795.44
$ws.Range("J99").Value = 16259.091
$ws.Range("K99").Value = 6795.44
$ws.Range("L99").Value = 16259.091
$ws.Range("M99").Value = -5297.44
$ws.Range("N99").Value = -19255.091
$ws.Range("H126").Value = 9687.111000000001
$ws.Range("I126").Value = 6795.44
$ws.Range("J126").Value = 16259.091
$ws.Range("K126").Value = 20386.32
$ws.Range("L126").Value = 48777.273
$ws.Range("M126").Value = -17916.32
$ws.Range("N126").Value = -53717.273
$ws.Range("H132").Value = 2855.65
$ws.Range("I132").Value = 2358.1765
$ws.Range("J132").Value = 5674.6665
$ws.Range("K132").Value = 7074.529500000001
$ws.Range("L132").Value = 17023.9995
$ws.Range("M132").Value = -4544.529500000001
$ws.Range("N132").Value = -22083.9995
$ws.Range("H141").Value = 82500
$ws.Range("J141").Value = 88333.336
$ws.Range("L141").Value = 88333.336
$ws.Range("N141").Value = -98693.336

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H107").Value = 1570.6666
$ws.Range("J107").Value = 1932.6666
$ws.Range("L107").Value = 5797.9998
$ws.Range("N107").Value = -9637.9998
$ws.Range("H137").Value = 4012.7144
$ws.Range("I137").Value = 4730
$ws.Range("J137").Value = 3614.2222
$ws.Range("K137").Value = 14190
$ws.Range("L137").Value = 10842.6666
$ws.Range("M137").Value = -9090
$ws.Range("N137").Value = -21042.6666
$ws.Range("H139").Value = 4752.3076
$ws.Range("I139").Value = 1864.8889
$ws.Range("K139").Value = 5594.6667
$ws.Range("M139").Value = -454.6666999999998

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 3516.7144
$ws.Range("J43").Value = 6339.6665
$ws.Range("L43").Value = 6339.6665
$ws.Range("N43").Value = -6641.6665
$ws.Range("H57").Value = 37499.5
$ws.Range("J57").Value = 59999
$ws.Range("L57").Value = 59999
$ws.Range("N57").Value = -61639
$ws.Range("H97").Value = 1176.4138
$ws.Range("I97").Value = 1229.6666
$ws.Range("K97").Value = 1229.6666
$ws.Range("M97").Value = -733.6666
$ws.Range("H99").Value = 7650.846
$ws.Range("I99").Value = 4087.3635
$ws.Range("K99").Value = 4087.3635
$ws.Range("M99").Value = -1841.3635
$ws.Range("H132").Value = 4513.1333
$ws.Range("I132").Value = 3985.3076
$ws.Range("J132").Value = 7944
$ws.Range("K132").Value = 11955.9228
$ws.Range("L132").Value = 23832
$ws.Range("M132").Value = -9425.9228
$ws.Range("N132").Value = -28892

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1575.5714
$ws.Range("I22").Value = 682.5
$ws.Range("J22").Value = 2766.3333
$ws.Range("K22").Value = 682.5
$ws.Range("L22").Value = 2766.3333
$ws.Range("M22").Value = -387.5
$ws.Range("N22").Value = -3356.3333
$ws.Range("H27").Value = 1575.5714
$ws.Range("I27").Value = 682.5
$ws.Range("J27").Value = 2766.3333
$ws.Range("K27").Value = 682.5
$ws.Range("L27").Value = 2766.3333
$ws.Range("M27").Value = -575.5
$ws.Range("N27").Value = -2980.3333

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7627.875
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 7627.875
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 7627.875
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -8875.875
$ws.Range("H65").Value = 7627.875
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 7627.875
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 38139.375
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -44379.375
